# ParameterInputsTrout.xlsx - add the new "Min/Resp Parameters" section
# (Hilary and Paul's mineralization routine) below the existing SWGW
# parameters block, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20: new section header, mirrors the existing "#...Parameters" rows
$ws.Range("A20").Value = "#Min/Resp Parameters"

# Row 21: the new parameter itself - name / value / units
$ws.Range("A21").Value = "DOC_miner_const"
$ws.Range("B21").Value = 0.005
$ws.Range("C21").Value = "1/days"

# Reflect the scrolled viewport + new selection that was saved with the
# workbook (best effort - selecting the next empty row below the data).
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A22").Select()
